$wb = $excel.ActiveWorkbook

# --- dimLanguage sheet: add "Tools" column and update German proficiency ---
$ws = $wb.Worksheets.Item("dimLanguage")

# New header cell E1 = "Tools" (copy header style from D1 so it matches s="1")
$ws.Range("E1").Value = "Tools"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New "Tools" values for the technical-skill rows
$ws.Range("E4").Value = "DAX Queries, Data Modeling, Dashboard Creation, Data Cleaning and Preparation"
$ws.Range("E5").Value = "Pandas, Numpy, SQLalchemy"
$ws.Range("E6").Value = "Advanced SQL Queries"
$ws.Range("E7").Value = "Advanced VBA"
$ws.Range("E8").Value = "Scikit"
$ws.Range("E9").Value = "Functional roles"
$ws.Range("E10").Value = "Functional roles"

# Update German proficiency from "C2" to "Goethe Institut - C2"
$ws.Range("C3").Value = "Goethe Institut - C2"

# Make dimLanguage the active/selected sheet with C4 selected (was fctLifePhases before)
$ws.Activate()
$null = $ws.Range("C4").Select()
